$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new rows of personal expenditure / data source links
$ws.Range("D11").Value = "https://beta.bls.gov/dataViewer/view/timeseries/LNS14000000;jsessionid=B67E9E4C2185836C45B62A6FC973C7EB"
$ws.Range("C11").Value = "BLS Time Series by Month for Unemployment 2008-2021"

$ws.Range("D12").Value = "https://www.esrl.noaa.gov/gmd/ccgg/trends/data.html"
$ws.Range("C12").Value = "NOAA CO2 emissions"

$ws.Range("C11:C12").Borders.Item(7).LineStyle = 1

$ws.Range("D19").Select()
